# SwaadSutra_Consolidated_2026-01-13.xlsx
# New order (Order ID 10 - "Yuu", Wheat Chapati x1) placed at 2026-01-13 19:12,
# inserted as the newest row at the top of the "All Orders" log (row 2),
# pushing every existing order row down by one. "Daily Summary" totals are
# bumped to reflect the extra order (+1 order, +15 revenue/pending).

$wb = $excel.ActiveWorkbook

# ---- "All Orders" sheet ------------------------------------------------
$ws = $wb.Worksheets.Item("All Orders")

# Push existing data rows (2..10) down to (3..11) and open up a fresh row 2
# for the newest order.
$ws.Rows("2:2").Insert()

# Columns E (Phone) and J (Collection Date) hold digit-only / date-look
# text in this log (phone numbers, ISO dates). Force them to Text first so
# Excel doesn't reinterpret the literal strings as a number / date serial -
# matching every other row in the sheet, where these columns are plain text.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("J2").NumberFormat = "@"

$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "2026-01-13 19:12"
$ws.Range("C2").Value = "Sagar Borse"
$ws.Range("D2").Value = "Yuu"
$ws.Range("E2").Value = "7588930329"
$ws.Range("F2").Value = "Wheat Chapati x1"
$ws.Range("G2").Value = 15
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"
$ws.Range("J2").Value = "2026-01-15"
$ws.Range("K2").Value = "02:42"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""

# ---- "Daily Summary" sheet ---------------------------------------------
$ws2 = $wb.Worksheets.Item("Daily Summary")

$ws2.Range("B2").Value = 10    # Total Orders: 9 -> 10
$ws2.Range("E2").Value = 245   # Revenue: 230 -> 245 (+15 for the new order)
$ws2.Range("G2").Value = 245   # Pending: 230 -> 245 (+15 for the new order)
